$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.894.11"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.066.05"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.16"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.06"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.30"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.593.45"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.882.28"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.066.67"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.10"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.61"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.501"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.45"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0903"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  +5.92%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.73"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.04"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.02"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.98"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0679"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.108.45"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.47"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.655"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.262.55"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0259"
$ws.Range("E45").Value = "  +7.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.83"
$ws.Range("E46").Value = "  +5.99%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.941"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.94"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.740"
$ws.Range("E50").Value = "  +9.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "256.42"
$ws.Range("E51").Value = "  +11.70%  "
